$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value2 = 0
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 0
$ws.Range("K21").Value2 = 0
$ws.Range("N21").Value2 = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
# Row 23
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 0
$ws.Range("N23").Value2 = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").ClearContents()
# Row 98
$ws.Range("H98").Value2 = 52083904
$ws.Range("I98").Value2 = 56818668
$ws.Range("J98").Value2 = 1500
$ws.Range("K98").Value2 = 56818668
$ws.Range("L98").Value2 = 1500
$ws.Range("M98").Value2 = -56817170
$ws.Range("N98").Value2 = -4496
# Row 122
$ws.Range("H122").Value2 = 52083904
$ws.Range("I122").Value2 = 56818668
$ws.Range("J122").Value2 = 1500
$ws.Range("K122").Value2 = 170456004
$ws.Range("L122").Value2 = 4500
$ws.Range("M122").Value2 = -170453554
$ws.Range("N122").Value2 = -9400
# Row 137
$ws.Range("H137").Value2 = 1347.9535
$ws.Range("I137").Value2 = 1289.64
$ws.Range("K137").Value2 = 3868.92
$ws.Range("M137").Value2 = -1318.92
# Row 138
$ws.Range("H138").Value2 = 1498.3334
$ws.Range("I138").Value2 = 867.2093
$ws.Range("J138").Value2 = 3436.7856
$ws.Range("K138").Value2 = 2601.6279
$ws.Range("L138").Value2 = 10310.3568
$ws.Range("M138").Value2 = 2538.3721
$ws.Range("N138").Value2 = -20590.3568
# Row 141
$ws.Range("H141").Value2 = 2168.4478
$ws.Range("I141").Value2 = 966.18335
$ws.Range("J141").Value2 = 12473.571
$ws.Range("K141").Value2 = 2898.55005
$ws.Range("L141").Value2 = 37420.713
$ws.Range("M141").Value2 = 2281.44995
$ws.Range("N141").Value2 = -47780.713

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 1694.1099
$ws.Range("I32").Value2 = 1613.0588
$ws.Range("J32").Value2 = 2842.3333
$ws.Range("K32").Value2 = 1613.0588
$ws.Range("L32").Value2 = 2842.3333
$ws.Range("M32").Value2 = -1326.0588
$ws.Range("N32").Value2 = -3416.3333
# Row 74
$ws.Range("H74").Value2 = 950.4706
$ws.Range("I74").Value2 = 969.5111000000001
$ws.Range("J74").Value2 = 807.6667
$ws.Range("K74").Value2 = 969.5111000000001
$ws.Range("L74").Value2 = 807.6667
$ws.Range("M74").Value2 = -95.51110000000006
$ws.Range("N74").Value2 = -2555.6667
# Row 77
$ws.Range("H77").Value2 = 950.4706
$ws.Range("I77").Value2 = 969.5111000000001
$ws.Range("J77").Value2 = 807.6667
$ws.Range("K77").Value2 = 4847.5555
$ws.Range("L77").Value2 = 4038.3335
$ws.Range("M77").Value2 = -479.5555000000004
$ws.Range("N77").Value2 = -12774.3335
# Row 132
$ws.Range("H132").Value2 = 19252336
$ws.Range("I132").Value2 = 23810290
$ws.Range("J132").Value2 = 4526642.5
$ws.Range("K132").Value2 = 71430870
$ws.Range("L132").Value2 = 13579927.5
$ws.Range("M132").Value2 = -71428340
$ws.Range("N132").Value2 = -13584987.5

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value2 = 4449584
$ws.Range("I134").Value2 = 1000.0833
$ws.Range("J134").Value2 = 8555970
$ws.Range("K134").Value2 = 3000.2499
$ws.Range("L134").Value2 = 25667910
$ws.Range("M134").Value2 = -465.2498999999998
$ws.Range("N134").Value2 = -25672980
# Row 139
$ws.Range("H139").Value2 = 46864.445
$ws.Range("I139").Value2 = 44990
$ws.Range("J139").Value2 = 47400
$ws.Range("K139").Value2 = 44990
$ws.Range("L139").Value2 = 47400
$ws.Range("M139").Value2 = -39850
$ws.Range("N139").Value2 = -57680

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 1976.7084
$ws.Range("I31").Value2 = 1583.3334
$ws.Range("J31").Value2 = 2107.8333
$ws.Range("K31").Value2 = 1583.3334
$ws.Range("L31").Value2 = 2107.8333
$ws.Range("M31").Value2 = -1288.3334
$ws.Range("N31").Value2 = -2697.8333
# Row 34
$ws.Range("H34").Value2 = 1976.7084
$ws.Range("I34").Value2 = 1583.3334
$ws.Range("J34").Value2 = 2107.8333
$ws.Range("K34").Value2 = 1583.3334
$ws.Range("L34").Value2 = 2107.8333
$ws.Range("M34").Value2 = -1381.3334
$ws.Range("N34").Value2 = -2511.8333
# Row 134
$ws.Range("H134").Value2 = 1408.1904
$ws.Range("I134").Value2 = 1271
$ws.Range("K134").Value2 = 3813
$ws.Range("M134").Value2 = -1278

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 27780676
$ws.Range("J5").Value2 = 27784086
$ws.Range("L5").Value2 = 83352258
$ws.Range("N5").Value2 = -83352482
# Row 98
$ws.Range("H98").Value2 = 918.82355
$ws.Range("I98").Value2 = 750
$ws.Range("J98").Value2 = 929.375
$ws.Range("K98").Value2 = 2250
$ws.Range("L98").Value2 = 2788.125
$ws.Range("M98").Value2 = -752
$ws.Range("N98").Value2 = -5784.125
# Row 107
$ws.Range("H107").Value2 = 23706410
$ws.Range("I107").Value2 = 153.6875
$ws.Range("J107").Value2 = 36785724
$ws.Range("K107").Value2 = 461.0625
$ws.Range("L107").Value2 = 110357172
$ws.Range("M107").Value2 = 1458.9375
$ws.Range("N107").Value2 = -110361012
# Row 113
$ws.Range("H113").Value2 = 19074468
$ws.Range("I113").Value2 = 7576158.5
$ws.Range("J113").Value2 = 25731384
$ws.Range("K113").Value2 = 22728475.5
$ws.Range("L113").Value2 = 77194152
$ws.Range("M113").Value2 = -22726305.5
$ws.Range("N113").Value2 = -77198492
# Row 135
$ws.Range("H135").Value2 = 27780676
$ws.Range("J135").Value2 = 27784086
$ws.Range("L135").Value2 = 250056774
$ws.Range("N135").Value2 = -250061844

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value2 = 3090.0908
$ws.Range("J126").Value2 = 3197.9
$ws.Range("L126").Value2 = 9593.700000000001
$ws.Range("N126").Value2 = -14533.7
# Row 132
$ws.Range("H132").Value2 = 5526.439
$ws.Range("I132").Value2 = 3516.4546
$ws.Range("J132").Value2 = 13817.625
$ws.Range("K132").Value2 = 10549.3638
$ws.Range("L132").Value2 = 41452.875
$ws.Range("M132").Value2 = -8019.363799999999
$ws.Range("N132").Value2 = -46512.875

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 1508.2142
$ws.Range("I7").Value2 = 1426.3334
$ws.Range("J7").Value2 = 1999.5
$ws.Range("K7").Value2 = 1426.3334
$ws.Range("L7").Value2 = 1999.5
$ws.Range("M7").Value2 = -1314.3334
$ws.Range("N7").Value2 = -2223.5
# Row 126
$ws.Range("H126").Value2 = 1508.2142
$ws.Range("I126").Value2 = 1426.3334
$ws.Range("J126").Value2 = 1999.5
$ws.Range("K126").Value2 = 4279.0002
$ws.Range("L126").Value2 = 5998.5
$ws.Range("M126").Value2 = -1809.0002
$ws.Range("N126").Value2 = -10938.5
